# PeerEvalTemplate_3Person.xlsx - typo / wording cleanup pass
#
# Matches the commit "Latest canvas course groups small cleanup in .r file /
# fixed typos in the blank templates": the rating-scale descriptions and the
# header instructions had a handful of small wording/typo fixes, and the
# active selection was left on E4 (the "Rating Descriptions" header cell)
# instead of A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CRLF = "`r`n"

# E5: KSA rating scale - "peform" -> "perform", "one elses" -> "one else's"
$e5 = @(
    "5: Demonstrates KSAs to do excellent work, acquires new KSA to help team, can perform any role on team if necessary",
    "4: Between 5 above and 3 below",
    "3: Demonstrates sufficient KSA to contribute to team, acquires KSAs to meet requirements, able to perform other tasks",
    "2: Between 3 above and 1 below",
    "1: Missing basic qualification, unable to develop KSAs to contribute to team, unable to perform any one else's duties"
) -join $CRLF
$ws.Range("E5").Value2 = $e5

# E7: Interacting with Teammates rating scale - "contributsions" -> "contributions"
$e7 = @(
    "5: Is interested in teammates ideas and contributions, makes sure everyone is informed, is encouraging, enthusiastic and asks for feedback/suggestions",
    "3: Listens and respects teammate contributions, communicates clearly, shares info, participates fully, reacts and responds to feedback/suggestions",
    "1: Interrupts, ignores, bosses, or makes fun, takes action without input, does not share, complains, makes excuses, does not interact, is defensive"
) -join $CRLF
$ws.Range("E7").Value2 = $e7

# E8: Keeping the Team on Track rating scale - "sucess" -> "success"
$e8 = @(
    "5: Monitors teams' progress, makes sure teammates are progressing, gives specific, timely, and constructive feedback",
    "3: Knows what everyone on the team should be doing and notices problems, alerts teammates and suggests solutions with success is threatened",
    "1: Unaware if team is meeting goals, does not pay attention to teammates progress, avoids discussing team problems even when obvious"
) -join $CRLF
$ws.Range("E8").Value2 = $e8

# E9: Expecting Quality rating scale - "Encouarges" -> "Encourages", "responsiblities" -> "responsibilities"
$e9 = @(
    "5: Motivates team to do excellent work, cares about excellent work even without reward, believes in team's ability to do excellent work",
    "3: Encourages good work to meet requirements, believes team can meet its responsibilities",
    "1: Satisfied even if not all requirements are met,  avoids work, doubts team can meet requirements"
) -join $CRLF
$ws.Range("E9").Value2 = $e9

# E4: "Category" header note - drop the duplicated word "each"
$ws.Range("E4").Value2 = "Rating Descriptions (provide whole number ratings (5, 4, 3, 2, or 1) in columns for each member including yourself)"

# Leave the cursor/selection on the header cell that was edited, matching the
# saved workbook's <selection activeCell="E4" sqref="E4"/>.
$ws.Range("E4").Select()
